$d = $word.ActiveDocument

# Helper: standard wildcard-off literal find & replace across whole document
function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# --- Paragraph 1, first big run ---
Replace-Text "studerende der er motiveret for at lærer mere inden for programmerings verden. På nuværende tidspunkt er jeg på mit 3. semester på kandidaten. Jeg søger et job hvor jeg kan anvende min opnået viden til at programmere løsninger." "studerende, der er motiveret for, at lære mere indenfor programmeringsverdenen. På nuværende tidspunkt er jeg på mit 3. semester på kandidaten. Jeg søger et job, hvor jeg kan anvende min opnåede viden til at programmere løsninger."

# --- Paragraph 1, "Herudover" sentence ---
Replace-Text "Herudover har jeg designet forskellige frontends til applikationer så jeg er bekendt med HTML" "Derudover har jeg designet forskellige frontends til applikationer, så jeg er bekendt med HTML"

# --- Paragraph 2, "patterns ... den objekt orienteret" ---
Replace-Text "patterns såsom MVC og har specielt viden indenfor den objekt orienteret tilgang." "patterns, såsom MVC og har specielt viden indenfor den objektorienterede tilgang."

# --- Paragraph 4, closing paragraph ---
Replace-Text "Jeg håber på at anvende ovenstående til at udfylde en organisations opgaver og mål i fremtiden." "Jeg håber på at anvende mine kompetencer til at gennemføre en organisations opgaver og mål i fremtiden."
